# Updated cryptos list values (Price / Volume(1h)) per the target diff.
# Cells hold text that LOOKS numeric (multi-dot prices, "+x.xx%  " strings),
# so we force Text format before writing and restore the default "Normal"
# style afterwards (matches original cells, which carry no style override).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 4; Value = "27.377.72" }
    @{ Row = 2; Col = 5; Value = "  +2.37%  " }
    @{ Row = 3; Col = 4; Value = "1.662.80" }
    @{ Row = 3; Col = 5; Value = "  +1.34%  " }
    @{ Row = 4; Col = 5; Value = "  -0.48%  " }
    @{ Row = 5; Col = 4; Value = "220.45" }
    @{ Row = 5; Col = 5; Value = "  +1.21%  " }
    @{ Row = 6; Col = 5; Value = "  +0.78%  " }
    @{ Row = 7; Col = 5; Value = "  -0.51%  " }
    @{ Row = 8; Col = 5; Value = "  +1.36%  " }
    @{ Row = 9; Col = 5; Value = "  +0.29%  " }
    @{ Row = 10; Col = 4; Value = "20.14" }
    @{ Row = 10; Col = 5; Value = "  +5.40%  " }
    @{ Row = 11; Col = 5; Value = "  +0.77%  " }
    @{ Row = 12; Col = 4; Value = "1.894.52" }
    @{ Row = 12; Col = 5; Value = "  +1.30%  " }
    @{ Row = 13; Col = 4; Value = "1.669.70" }
    @{ Row = 13; Col = 5; Value = "  +1.65%  " }
    @{ Row = 14; Col = 5; Value = "  +0.86%  " }
    @{ Row = 15; Col = 4; Value = "0.534" }
    @{ Row = 15; Col = 5; Value = "  +1.21%  " }
    @{ Row = 16; Col = 4; Value = "67.18" }
    @{ Row = 16; Col = 5; Value = "  +3.86%  " }
    @{ Row = 17; Col = 4; Value = "27.373.63" }
    @{ Row = 17; Col = 5; Value = "  +2.39%  " }
    @{ Row = 18; Col = 5; Value = "  +0.58%  " }
    @{ Row = 19; Col = 4; Value = "222.83" }
    @{ Row = 19; Col = 5; Value = "  +4.03%  " }
    @{ Row = 20; Col = 5; Value = "  -0.49%  " }
    @{ Row = 21; Col = 4; Value = "6.76" }
    @{ Row = 21; Col = 5; Value = "  +8.68%  " }
    @{ Row = 22; Col = 5; Value = "  +1.69%  " }
    @{ Row = 23; Col = 4; Value = "2.50" }
    @{ Row = 23; Col = 5; Value = "  +5.68%  " }
    @{ Row = 24; Col = 4; Value = "9.30" }
    @{ Row = 24; Col = 5; Value = "  +0.31%  " }
    @{ Row = 25; Col = 4; Value = "147.06" }
    @{ Row = 25; Col = 5; Value = "  +1.01%  " }
    @{ Row = 26; Col = 5; Value = "  -0.64%  " }
    @{ Row = 27; Col = 4; Value = "7.44" }
    @{ Row = 27; Col = 5; Value = "  +3.93%  " }
    @{ Row = 28; Col = 4; Value = "0.120" }
    @{ Row = 28; Col = 5; Value = "  +1.20%  " }
    @{ Row = 29; Col = 5; Value = "  +2.60%  " }
    @{ Row = 30; Col = 4; Value = "0.0514" }
    @{ Row = 30; Col = 5; Value = "  +1.03%  " }
    @{ Row = 31; Col = 5; Value = "  +0.96%  " }
    @{ Row = 32; Col = 4; Value = "3.40" }
    @{ Row = 32; Col = 5; Value = "  +0.11%  " }
    @{ Row = 33; Col = 5; Value = "  +0.04%  " }
    @{ Row = 34; Col = 5; Value = "  +2.34%  " }
    @{ Row = 35; Col = 4; Value = "1.265.96" }
    @{ Row = 35; Col = 5; Value = "  -1.65%  " }
    @{ Row = 36; Col = 5; Value = "  +0.95%  " }
    @{ Row = 37; Col = 4; Value = "0.0179" }
    @{ Row = 37; Col = 5; Value = "  +0.83%  " }
    @{ Row = 38; Col = 4; Value = "0.539" }
    @{ Row = 38; Col = 5; Value = "  +0.28%  " }
    @{ Row = 39; Col = 4; Value = "0.839" }
    @{ Row = 39; Col = 5; Value = "  +2.71%  " }
    @{ Row = 40; Col = 5; Value = "  -0.44%  " }
    @{ Row = 41; Col = 4; Value = "0.819" }
    @{ Row = 41; Col = 5; Value = "  +1.63%  " }
    @{ Row = 42; Col = 4; Value = "5.40" }
    @{ Row = 42; Col = 5; Value = "  +2.24%  " }
    @{ Row = 43; Col = 4; Value = "1.806.56" }
    @{ Row = 43; Col = 5; Value = "  +1.52%  " }
    @{ Row = 44; Col = 5; Value = "  -4.08%  " }
    @{ Row = 45; Col = 4; Value = "61.87" }
    @{ Row = 45; Col = 5; Value = "  +1.38%  " }
    @{ Row = 46; Col = 4; Value = "92.37" }
    @{ Row = 46; Col = 5; Value = "  +0.80%  " }
    @{ Row = 47; Col = 5; Value = "  +1.54%  " }
    @{ Row = 48; Col = 5; Value = "  +0.19%  " }
    @{ Row = 49; Col = 4; Value = "0.0987" }
    @{ Row = 49; Col = 5; Value = "  +2.26%  " }
    @{ Row = 50; Col = 4; Value = "7.71" }
    @{ Row = 50; Col = 5; Value = "  +1.15%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
